# HP: Testing 2 testcases
# Adds a second set of testcase(2) result/error columns to the "Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# New header cells K10:O10 should carry the same bordered header style as
# the rest of row 10 (A10:J10 use style index 2). Copy the format from the
# existing J10 header cell instead of re-creating a style object, so the
# saved workbook reuses the existing cellXfs entry instead of minting a
# near-duplicate one.
$ws.Range("J10").Copy()
$ws.Range("K10:O10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(10, 12).Value = "Errors TC1"   # L10
$ws.Cells.Item(10, 13).Value = "Errors TC2"   # M10
$ws.Cells.Item(10, 14).Value = "Errors TC3"   # N10
$ws.Cells.Item(10, 15).Value = "Errors TC4"   # O10
# K10 stays blank (format only, no value).

# Testcase(2) result messages in column H, rows 11-14 (mirrors the existing
# Testcase(1) messages already present in column G).
$ws.Cells.Item(11, 8).Value = "RESULT: SUBMIT(1) TESTCASE(2) MSG:(TESTCASE#2:EQUAL`n)"
$ws.Cells.Item(12, 8).Value = "RESULT: SUBMIT(2) TESTCASE(2) MSG:(TESTCASE#2:NOT_EQUAL`n)"
$ws.Cells.Item(13, 8).Value = "RESULT: SUBMIT(3) TESTCASE(2) MSG:(TESTCASE#2:NOT_EQUAL`n)"
$ws.Cells.Item(14, 8).Value = "RESULT: SUBMIT(4) TESTCASE(2) MSG:(TESTCASE#2:NOT_EQUAL`n)"

# Error messages for submit 4 testcase 1 (moved from H14 to N14) and the new
# testcase 2 error (O14).
$ws.Cells.Item(14, 14).Value = "ERROR: SUBMIT(4) TESTCASE(1) MSG:(CLASS:siima.app.XSLTransformer ERROR:Syntax error in '/CATALOG/Plant[(Price>'6.60')] and [(Light='Sun')]'.)"
$ws.Cells.Item(14, 15).Value = "ERROR: SUBMIT(4) TESTCASE(2) MSG:(CLASS:siima.app.XSLTransformer ERROR:Syntax error in '/CATALOG/Plant[(Price>'6.60')] and [(Light='Sun')]'.)"

# Re-run autofit on the rows touched above so the embedded newlines in the
# messages don't leave a stale explicit row height behind.
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()

# New column widths for G:O (Excel auto-sized these when the content changed).
$ws.Columns.Item(7).ColumnWidth = 59.666666667
$ws.Columns.Item(8).ColumnWidth = 67.333333333
$ws.Columns.Item(9).ColumnWidth = 57.333333333
$ws.Columns.Item(10).ColumnWidth = 31.833333333
$ws.Columns.Item(12).ColumnWidth = 15.5
$ws.Columns.Item(13).ColumnWidth = 16.666666667
$ws.Columns.Item(14).ColumnWidth = 14.666666667
$ws.Columns.Item(15).ColumnWidth = 15.5

# Refresh the view to match the new selection / scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("L15").Select()
